# Auto-generated Excel COM-interop script to apply workbook value updates
$wb = $excel.ActiveWorkbook

# Sheet ALC (index 1), row 9
$ws = $wb.Worksheets.Item(1)
$ws.Range("H9").Value = 991.6667
$ws.Range("I9").Value = 1059.5
$ws.Range("K9").Value = 1059.5
$ws.Range("M9").Value = -890.5

# Sheet ALC (index 1), row 86
$ws = $wb.Worksheets.Item(1)
$ws.Range("H86").Value = 7947.1577
$ws.Range("I86").Value = 8554.666999999999
$ws.Range("J86").Value = 7400.4
$ws.Range("K86").Value = 8554.666999999999
$ws.Range("L86").Value = 7400.4
$ws.Range("M86").Value = -7431.666999999999
$ws.Range("N86").Value = -9646.4

# Sheet ALC (index 1), row 89
$ws = $wb.Worksheets.Item(1)
$ws.Range("H89").Value = 7947.1577
$ws.Range("I89").Value = 8554.666999999999
$ws.Range("J89").Value = 7400.4
$ws.Range("K89").Value = 42773.335
$ws.Range("L89").Value = 37002
$ws.Range("M89").Value = -37157.335
$ws.Range("N89").Value = -48234

# Sheet ALC (index 1), row 138
$ws = $wb.Worksheets.Item(1)
$ws.Range("H138").Value = 4915
$ws.Range("J138").Value = 4713.3
$ws.Range("L138").Value = 14139.9
$ws.Range("N138").Value = -24419.9

# Sheet ARM (index 2), row 57
$ws = $wb.Worksheets.Item(2)
$ws.Range("H57").Value = 14975
$ws.Range("I57").Value = 14975
$ws.Range("K57").Value = 14975
$ws.Range("M57").Value = -14491

# Sheet ARM (index 2), row 88
$ws = $wb.Worksheets.Item(2)
$ws.Range("H88").Value = 2328.1428
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2328.1428
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 2328.1428
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -3140.1428

# Sheet ARM (index 2), row 91
$ws = $wb.Worksheets.Item(2)
$ws.Range("H91").Value = 2328.1428
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2328.1428
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 2328.1428
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -5136.1428

# Sheet ARM (index 2), row 110
$ws = $wb.Worksheets.Item(2)
$ws.Range("H110").Value = 4106.5
$ws.Range("I110").Value = 4861.0645
$ws.Range("K110").Value = 4861.0645
$ws.Range("M110").Value = -2816.0645

# Sheet BSM (index 3), row 108
$ws = $wb.Worksheets.Item(3)
$ws.Range("H108").Value = 74999.5
$ws.Range("J108").Value = 74999.5
$ws.Range("L108").Value = 74999.5
$ws.Range("N108").Value = -82679.5

# Sheet BSM (index 3), row 109
$ws = $wb.Worksheets.Item(3)
$ws.Range("H109").Value = 58684
$ws.Range("J109").Value = 58684
$ws.Range("L109").Value = 58684
$ws.Range("N109").Value = -61458

# Sheet BSM (index 3), row 134
$ws = $wb.Worksheets.Item(3)
$ws.Range("H134").Value = 68510.766
$ws.Range("I134").Value = 107518
$ws.Range("J134").Value = 25602.8
$ws.Range("K134").Value = 322554
$ws.Range("L134").Value = 76808.39999999999
$ws.Range("M134").Value = -320019
$ws.Range("N134").Value = -81878.39999999999

# Sheet BSM (index 3), row 140
$ws = $wb.Worksheets.Item(3)
$ws.Range("H140").Value = 100000
$ws.Range("J140").Value = 100000
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360

# Sheet CRP (index 4), row 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("H4").Value = 264.2143
$ws.Range("J4").Value = 649.5
$ws.Range("L4").Value = 649.5
$ws.Range("N4").Value = -873.5

# Sheet CRP (index 4), row 16
$ws = $wb.Worksheets.Item(4)
$ws.Range("H16").Value = 4375.2915
$ws.Range("I16").Value = 2296.4546
$ws.Range("J16").Value = 6134.3076
$ws.Range("K16").Value = 2296.4546
$ws.Range("L16").Value = 6134.3076
$ws.Range("M16").Value = -2009.4546
$ws.Range("N16").Value = -6708.3076

# Sheet CRP (index 4), row 62
$ws = $wb.Worksheets.Item(4)
$ws.Range("H62").Value = 2274.75
$ws.Range("I62").Value = 1250
$ws.Range("K62").Value = 1250
$ws.Range("M62").Value = -626

# Sheet CRP (index 4), row 65
$ws = $wb.Worksheets.Item(4)
$ws.Range("H65").Value = 2274.75
$ws.Range("I65").Value = 1250
$ws.Range("K65").Value = 6250
$ws.Range("M65").Value = -3130

# Sheet CRP (index 4), row 107
$ws = $wb.Worksheets.Item(4)
$ws.Range("H107").Value = 1907.2307
$ws.Range("I107").Value = 1752.8823
$ws.Range("J107").Value = 2198.7778
$ws.Range("K107").Value = 1752.8823
$ws.Range("L107").Value = 2198.7778
$ws.Range("M107").Value = 167.1177
$ws.Range("N107").Value = -6038.7778

# Sheet CRP (index 4), row 110
$ws = $wb.Worksheets.Item(4)
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

# Sheet CRP (index 4), row 113
$ws = $wb.Worksheets.Item(4)
$ws.Range("H113").Value = 4375.2915
$ws.Range("I113").Value = 2296.4546
$ws.Range("J113").Value = 6134.3076
$ws.Range("K113").Value = 2296.4546
$ws.Range("L113").Value = 6134.3076
$ws.Range("M113").Value = -126.4546
$ws.Range("N113").Value = -10474.3076

# Sheet CUL (index 5), row 7
$ws = $wb.Worksheets.Item(5)
$ws.Range("H7").Value = 5242.8096
$ws.Range("I7").Value = 122.111115
$ws.Range("K7").Value = 366.333345
$ws.Range("M7").Value = -254.333345

# Sheet CUL (index 5), row 55
$ws = $wb.Worksheets.Item(5)
$ws.Range("H55").Value = 1487.5
$ws.Range("J55").Value = 1463
$ws.Range("L55").Value = 4389
$ws.Range("N55").Value = -4743

# Sheet CUL (index 5), row 80
$ws = $wb.Worksheets.Item(5)
$ws.Range("H80").Value = 15669.429
$ws.Range("I80").Value = 1800
$ws.Range("J80").Value = 17981
$ws.Range("K80").Value = 5400
$ws.Range("L80").Value = 53943
$ws.Range("M80").Value = -4464
$ws.Range("N80").Value = -55815

# Sheet CUL (index 5), row 81
$ws = $wb.Worksheets.Item(5)
$ws.Range("H81").Value = 8502000
$ws.Range("I81").Value = 3013
$ws.Range("J81").Value = 10201798
$ws.Range("K81").Value = 9039
$ws.Range("L81").Value = 30605394
$ws.Range("M81").Value = -7916
$ws.Range("N81").Value = -30607640

# Sheet CUL (index 5), row 83
$ws = $wb.Worksheets.Item(5)
$ws.Range("H83").Value = 15669.429
$ws.Range("I83").Value = 1800
$ws.Range("J83").Value = 17981
$ws.Range("K83").Value = 16200
$ws.Range("L83").Value = 161829
$ws.Range("M83").Value = -11520
$ws.Range("N83").Value = -171189

# Sheet CUL (index 5), row 84
$ws = $wb.Worksheets.Item(5)
$ws.Range("H84").Value = 8502000
$ws.Range("I84").Value = 3013
$ws.Range("J84").Value = 10201798
$ws.Range("K84").Value = 27117
$ws.Range("L84").Value = 91816182
$ws.Range("M84").Value = -21501
$ws.Range("N84").Value = -91827414

# Sheet CUL (index 5), row 131
$ws = $wb.Worksheets.Item(5)
$ws.Range("H131").Value = 1491.02
$ws.Range("I131").Value = 1354.3334
$ws.Range("J131").Value = 1495.2474
$ws.Range("K131").Value = 4063.0002
$ws.Range("L131").Value = 4485.7422
$ws.Range("M131").Value = 976.9998000000001
$ws.Range("N131").Value = -14565.7422

# Sheet CUL (index 5), row 139
$ws = $wb.Worksheets.Item(5)
$ws.Range("H139").Value = 9036.275
$ws.Range("I139").Value = 12542.462
$ws.Range("K139").Value = 37627.386
$ws.Range("M139").Value = -32487.386

# Sheet GSM (index 6), row 52
$ws = $wb.Worksheets.Item(6)
$ws.Range("H52").Value = 23777.525
$ws.Range("I52").Value = 20000
$ws.Range("J52").Value = 23879.621
$ws.Range("K52").Value = 20000
$ws.Range("L52").Value = 23879.621
$ws.Range("M52").Value = -19741
$ws.Range("N52").Value = -24397.621

# Sheet GSM (index 6), row 70
$ws = $wb.Worksheets.Item(6)
$ws.Range("H70").Value = 7422.727
$ws.Range("I70").Value = 5163.2856
$ws.Range("J70").Value = 11376.75
$ws.Range("K70").Value = 5163.2856
$ws.Range("L70").Value = 11376.75
$ws.Range("M70").Value = -4893.2856
$ws.Range("N70").Value = -11916.75

# Sheet GSM (index 6), row 73
$ws = $wb.Worksheets.Item(6)
$ws.Range("H73").Value = 7422.727
$ws.Range("I73").Value = 5163.2856
$ws.Range("J73").Value = 11376.75
$ws.Range("K73").Value = 5163.2856
$ws.Range("L73").Value = 11376.75
$ws.Range("M73").Value = -4227.2856
$ws.Range("N73").Value = -13248.75

# Sheet GSM (index 6), row 104
$ws = $wb.Worksheets.Item(6)
$ws.Range("H104").Value = 46835.5
$ws.Range("J104").Value = 46835.5
$ws.Range("L104").Value = 46835.5
$ws.Range("N104").Value = -53823.5

# Sheet GSM (index 6), row 126
$ws = $wb.Worksheets.Item(6)
$ws.Range("H126").Value = 7919.55
$ws.Range("I126").Value = 11506.333
$ws.Range("J126").Value = 6382.357
$ws.Range("K126").Value = 34518.999
$ws.Range("L126").Value = 19147.071
$ws.Range("M126").Value = -32048.999
$ws.Range("N126").Value = -24087.071

# Sheet GSM (index 6), row 132
$ws = $wb.Worksheets.Item(6)
$ws.Range("H132").Value = 960275.7
$ws.Range("I132").Value = 6733.1665
$ws.Range("J132").Value = 2867360.8
$ws.Range("K132").Value = 20199.4995
$ws.Range("L132").Value = 8602082.399999999
$ws.Range("M132").Value = -17669.4995
$ws.Range("N132").Value = -8607142.399999999

# Sheet GSM (index 6), row 137
$ws = $wb.Worksheets.Item(6)
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# Sheet GSM (index 6), row 139
$ws = $wb.Worksheets.Item(6)
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# Sheet LTW (index 7), row 46
$ws = $wb.Worksheets.Item(7)
$ws.Range("H46").Value = 900
$ws.Range("J46").Value = 900
$ws.Range("L46").Value = 900
$ws.Range("N46").Value = -1276

# Sheet LTW (index 7), row 107
$ws = $wb.Worksheets.Item(7)
$ws.Range("H107").Value = 7888.6665
$ws.Range("I107").Value = 7888.6665
$ws.Range("K107").Value = 7888.6665
$ws.Range("M107").Value = -5968.6665

# Sheet LTW (index 7), row 132
$ws = $wb.Worksheets.Item(7)
$ws.Range("H132").Value = 2410754
$ws.Range("I132").Value = 4240.6
$ws.Range("J132").Value = 7758562
$ws.Range("K132").Value = 12721.8
$ws.Range("L132").Value = 23275686
$ws.Range("M132").Value = -10191.8
$ws.Range("N132").Value = -23280746

# Sheet LTW (index 7), row 140
$ws = $wb.Worksheets.Item(7)
$ws.Range("H140").Value = 124729.75
$ws.Range("J140").Value = 124729.75
$ws.Range("L140").Value = 124729.75
$ws.Range("N140").Value = -135089.75

# Sheet WVR (index 8), row 81
$ws = $wb.Worksheets.Item(8)
$ws.Range("H81").Value = 1212.75
$ws.Range("I81").Value = 1212.75
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2425.5
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -1364.5
$ws.Range("N81").ClearContents()

# Sheet WVR (index 8), row 84
$ws = $wb.Worksheets.Item(8)
$ws.Range("H84").Value = 1212.75
$ws.Range("I84").Value = 1212.75
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 12127.5
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -6823.5
$ws.Range("N84").ClearContents()

# Sheet WVR (index 8), row 133
$ws = $wb.Worksheets.Item(8)
$ws.Range("H133").Value = 59366.25
$ws.Range("J133").Value = 59366.25
$ws.Range("L133").Value = 59366.25
$ws.Range("N133").Value = -69486.25

# Sheet WVR (index 8), row 138
$ws = $wb.Worksheets.Item(8)
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# Sheet WVR (index 8), row 140
$ws = $wb.Worksheets.Item(8)
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = 0
